$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.176.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.36%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.853.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'235.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.62%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4664"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.83%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2819"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.06412"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.92%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'18.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.84%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'97.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +14.60%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.857.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.08%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07543"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.26%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.982"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.47%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6387"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.12%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'294.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +21.13%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.171.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.000007378"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.07%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.096.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.050"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.41%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'164.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.095"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.50%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'19.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +7.34%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.932"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.19%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.1085"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.329"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.73%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.017"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.808"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04922"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7261"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.113"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.22%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.745"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.35%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.01920"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.668"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.64%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.8653"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.67%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.47%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'105.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.08%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.003"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.25%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.616"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.16%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.4054"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.89%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'65.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.93%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.90%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'8.981"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.51%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -2.17%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'34.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.55%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05557"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.3731"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.70%  "
$ws.Range("E51").Style = "Normal"
